$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 14.511828908386349
$ws.Range("C2").Value = 2.6473639617721005
$ws.Range("D2").Value = 1.1081228576872775
$ws.Range("E2").Value = 1.1522195782138169

$ws.Range("B3").Value = 5.2032961379966878
$ws.Range("C3").Value = 11.275118960341871
$ws.Range("D3").Value = 3.281647899050256
$ws.Range("E3").Value = -0.038992220035595437

$ws.Range("B1:E3").Select()
